$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -17.57306280842938
$ws.Range("C2").Value = 1.40042436625325
$ws.Range("D2").Value = -17.57306280842938
$ws.Range("E2").Value = -17.57306280842938
$ws.Range("F2").Value = -17.57306280842938
$ws.Range("G2").Value = -17.57306280842938
$ws.Range("H2").Value = -17.57306280842938
$ws.Range("I2").Value = -17.57306280842938
$ws.Range("J2").Value = -17.57306280842938
$ws.Range("K2").Value = -17.57306280842938
$ws.Range("B3").Value = -17.57306280842938
$ws.Range("C3").Value = -17.57306280842938
$ws.Range("D3").Value = -17.57306280842938
$ws.Range("E3").Value = -17.57306280842938
$ws.Range("F3").Value = -17.57306280842938
$ws.Range("G3").Value = -17.57306280842938
$ws.Range("H3").Value = -17.57306280842938
$ws.Range("I3").Value = 2.701175416779238
$ws.Range("J3").Value = -17.57306280842938
$ws.Range("K3").Value = -17.57306280842938
$ws.Range("B4").Value = -17.57306280842938
$ws.Range("C4").Value = 1.697782276492763
$ws.Range("D4").Value = 4.321921066103624
$ws.Range("E4").Value = -17.57306280842938
$ws.Range("F4").Value = 3.287145930299787
$ws.Range("G4").Value = -17.57306280842938
$ws.Range("H4").Value = 1.71654411806813
$ws.Range("I4").Value = -17.57306280842938
$ws.Range("J4").Value = 2.15776054253406
$ws.Range("K4").Value = -17.57306280842938
$ws.Range("B5").Value = -17.57306280842938
$ws.Range("C5").Value = 2.167872457514683
$ws.Range("D5").Value = -17.57306280842938
$ws.Range("E5").Value = -17.57306280842938
$ws.Range("F5").Value = -17.57306280842938
$ws.Range("G5").Value = 2.837054140732106
$ws.Range("H5").Value = -17.57306280842938
$ws.Range("I5").Value = -17.57306280842938
$ws.Range("J5").Value = -17.57306280842938
$ws.Range("K5").Value = -17.57306280842938
$ws.Range("B6").Value = -17.57306280842938
$ws.Range("C6").Value = -17.57306280842938
$ws.Range("D6").Value = -17.57306280842938
$ws.Range("E6").Value = -17.57306280842938
$ws.Range("F6").Value = -17.57306280842938
$ws.Range("G6").Value = -17.57306280842938
$ws.Range("H6").Value = -17.57306280842938
$ws.Range("I6").Value = -17.57306280842938
$ws.Range("J6").Value = -17.57306280842938
$ws.Range("K6").Value = -17.57306280842938
$ws.Range("B7").Value = 2.566199994092813
$ws.Range("C7").Value = -17.57306280842938
$ws.Range("D7").Value = -17.57306280842938
$ws.Range("E7").Value = -17.57306280842938
$ws.Range("F7").Value = -17.57306280842938
$ws.Range("G7").Value = -17.57306280842938
$ws.Range("H7").Value = -17.57306280842938
$ws.Range("I7").Value = -17.57306280842938
$ws.Range("J7").Value = -17.57306280842938
$ws.Range("K7").Value = -17.57306280842938
$ws.Range("B8").Value = -17.57306280842938
$ws.Range("C8").Value = -17.57306280842938
$ws.Range("D8").Value = -17.57306280842938
$ws.Range("E8").Value = 1.30746835805656
$ws.Range("F8").Value = -17.57306280842938
$ws.Range("G8").Value = -17.57306280842938
$ws.Range("H8").Value = -17.57306280842938
$ws.Range("I8").Value = -17.57306280842938
$ws.Range("J8").Value = -17.57306280842938
$ws.Range("K8").Value = -17.57306280842938
$ws.Range("B9").Value = 3.815312435697138
$ws.Range("C9").Value = -17.57306280842938
$ws.Range("D9").Value = -17.57306280842938
$ws.Range("E9").Value = -17.57306280842938
$ws.Range("F9").Value = -17.57306280842938
$ws.Range("G9").Value = -17.57306280842938
$ws.Range("H9").Value = -17.57306280842938
$ws.Range("I9").Value = -17.57306280842938
$ws.Range("J9").Value = -17.57306280842938
$ws.Range("K9").Value = -17.57306280842938
$ws.Range("B10").Value = -17.57306280842938
$ws.Range("C10").Value = -17.57306280842938
$ws.Range("D10").Value = -17.57306280842938
$ws.Range("E10").Value = -17.57306280842938
$ws.Range("F10").Value = -17.57306280842938
$ws.Range("G10").Value = -17.57306280842938
$ws.Range("H10").Value = -17.57306280842938
$ws.Range("I10").Value = 1.333713310985068
$ws.Range("J10").Value = -17.57306280842938
$ws.Range("K10").Value = 1.92862084876187
$ws.Range("B11").Value = -17.57306280842938
$ws.Range("C11").Value = -17.57306280842938
$ws.Range("D11").Value = -17.57306280842938
$ws.Range("E11").Value = 3.1324354336428
$ws.Range("F11").Value = -17.57306280842938
$ws.Range("G11").Value = 2.797161636175786
$ws.Range("H11").Value = -17.57306280842938
$ws.Range("I11").Value = -17.57306280842938
$ws.Range("J11").Value = -17.57306280842938
$ws.Range("K11").Value = 1.849406492464426
$ws.Range("B12").Value = -17.57306280842938
$ws.Range("C12").Value = -17.57306280842938
$ws.Range("D12").Value = -17.57306280842938
$ws.Range("E12").Value = -17.57306280842938
$ws.Range("F12").Value = -17.57306280842938
$ws.Range("G12").Value = -17.57306280842938
$ws.Range("H12").Value = -17.57306280842938
$ws.Range("I12").Value = -17.57306280842938
$ws.Range("J12").Value = -17.57306280842938
$ws.Range("K12").Value = -17.57306280842938
$ws.Range("B13").Value = -17.57306280842938
$ws.Range("C13").Value = -17.57306280842938
$ws.Range("D13").Value = -17.57306280842938
$ws.Range("E13").Value = 2.470766847267053
$ws.Range("F13").Value = -17.57306280842938
$ws.Range("G13").Value = -17.57306280842938
$ws.Range("H13").Value = -17.57306280842938
$ws.Range("I13").Value = -17.57306280842938
$ws.Range("J13").Value = 2.053075768768371
$ws.Range("K13").Value = 1.886218708449856
$ws.Range("B14").Value = -17.57306280842938
$ws.Range("C14").Value = -17.57306280842938
$ws.Range("D14").Value = -17.57306280842938
$ws.Range("E14").Value = -17.57306280842938
$ws.Range("F14").Value = -17.57306280842938
$ws.Range("G14").Value = -17.57306280842938
$ws.Range("H14").Value = -17.57306280842938
$ws.Range("I14").Value = -17.57306280842938
$ws.Range("J14").Value = -17.57306280842938
$ws.Range("K14").Value = 2.088021117588551
$ws.Range("B15").Value = -17.57306280842938
$ws.Range("C15").Value = -17.57306280842938
$ws.Range("D15").Value = -17.57306280842938
$ws.Range("E15").Value = -17.57306280842938
$ws.Range("F15").Value = -17.57306280842938
$ws.Range("G15").Value = -17.57306280842938
$ws.Range("H15").Value = -17.57306280842938
$ws.Range("I15").Value = -17.57306280842938
$ws.Range("J15").Value = -17.57306280842938
$ws.Range("K15").Value = -17.57306280842938
$ws.Range("B16").Value = -17.57306280842938
$ws.Range("C16").Value = -17.57306280842938
$ws.Range("D16").Value = -17.57306280842938
$ws.Range("E16").Value = -17.57306280842938
$ws.Range("F16").Value = -17.57306280842938
$ws.Range("G16").Value = -17.57306280842938
$ws.Range("H16").Value = -17.57306280842938
$ws.Range("I16").Value = -17.57306280842938
$ws.Range("J16").Value = 2.214815844397307
$ws.Range("K16").Value = -17.57306280842938
$ws.Range("B17").Value = -17.57306280842938
$ws.Range("C17").Value = 1.575763779472948
$ws.Range("D17").Value = -17.57306280842938
$ws.Range("E17").Value = -17.57306280842938
$ws.Range("F17").Value = -17.57306280842938
$ws.Range("G17").Value = -17.57306280842938
$ws.Range("H17").Value = 1.260559267787074
$ws.Range("I17").Value = 1.920552270620836
$ws.Range("J17").Value = 1.936296622190386
$ws.Range("K17").Value = -17.57306280842938
$ws.Range("B18").Value = -17.57306280842938
$ws.Range("C18").Value = -17.57306280842938
$ws.Range("D18").Value = -17.57306280842938
$ws.Range("E18").Value = -17.57306280842938
$ws.Range("F18").Value = -17.57306280842938
$ws.Range("G18").Value = -17.57306280842938
$ws.Range("H18").Value = 1.522779733434993
$ws.Range("I18").Value = 1.155272940681302
$ws.Range("J18").Value = 1.545160975271361
$ws.Range("K18").Value = -17.57306280842938
$ws.Range("B19").Value = -17.57306280842938
$ws.Range("C19").Value = -17.57306280842938
$ws.Range("D19").Value = -17.57306280842938
$ws.Range("E19").Value = -17.57306280842938
$ws.Range("F19").Value = -17.57306280842938
$ws.Range("G19").Value = -17.57306280842938
$ws.Range("H19").Value = 1.591316304628171
$ws.Range("I19").Value = 1.64264318073753
$ws.Range("J19").Value = -17.57306280842938
$ws.Range("K19").Value = -17.57306280842938
$ws.Range("B20").Value = -17.57306280842938
$ws.Range("C20").Value = 1.799396317598478
$ws.Range("D20").Value = -17.57306280842938
$ws.Range("E20").Value = -17.57306280842938
$ws.Range("F20").Value = 3.355878386288345
$ws.Range("G20").Value = -17.57306280842938
$ws.Range("H20").Value = 1.979525618493089
$ws.Range("I20").Value = 0.8803243902712536
$ws.Range("J20").Value = -17.57306280842938
$ws.Range("K20").Value = 2.214497052340562
$ws.Range("B21").Value = -17.57306280842938
$ws.Range("C21").Value = 1.660306430713241
$ws.Range("D21").Value = -17.57306280842938
$ws.Range("E21").Value = 1.683618118065658
$ws.Range("F21").Value = -17.57306280842938
$ws.Range("G21").Value = 2.561595382078276
$ws.Range("H21").Value = 2.165829355519956
$ws.Range("I21").Value = -17.57306280842938
$ws.Range("J21").Value = -17.57306280842938
$ws.Range("K21").Value = -17.57306280842938
